# Append the new daily allocation row (11/21/2025) to Sheet1.
# Column A stores the date as literal text (matching the existing rows,
# which are plain text strings rather than Excel date serials), so we
# force a text number format before assigning the value, then restore
# the cell to the default "Normal" style so no stray formatting differs
# from the other data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 81

$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "11/21/2025"
$dateCell.Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = 0.2092672469831965
$ws.Cells.Item($newRow, 3).Value = 0.7907327530168035
